$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert the last few sprint-3 contribution entries (rows 15-17),
# keeping the existing cell style but clearing the values.
$ws.Range("B15:D17").ClearContents()

# Update the current selection to match the reverted state.
$ws.Range("D14").Select()
